# Bulk market-price data refresh across all 8 profession sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates the price/profit columns (H:N) -- currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -- for
# specific leve rows with freshly pulled Universalis market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 126.9
$ws.Range("I33").Value = 85.44444
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 85.44444
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 143.55556
$ws.Range("N33").Value = -958
# Row 111
$ws.Range("H111").Value = 413
$ws.Range("I111").Value = 437.72726
$ws.Range("J111").Value = 277
$ws.Range("K111").Value = 1313.18178
$ws.Range("L111").Value = 831
$ws.Range("M111").Value = 1753.81822
$ws.Range("N111").Value = -6965
# Row 138
$ws.Range("H138").Value = 4371.525
$ws.Range("J138").Value = 5313.0586
$ws.Range("L138").Value = 15939.1758
$ws.Range("N138").Value = -26219.1758

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3846.6667
$ws.Range("I45").Value = 3731.4285
$ws.Range("K45").Value = 3731.4285
$ws.Range("M45").Value = -3354.4285
# Row 48
$ws.Range("H48").Value = 60000
$ws.Range("J48").Value = 60000
$ws.Range("L48").Value = 60000
$ws.Range("N48").Value = -60768
# Row 74
$ws.Range("H74").Value = 1120.2778
$ws.Range("I74").Value = 1004.4
$ws.Range("K74").Value = 1004.4
$ws.Range("M74").Value = -130.4
# Row 77
$ws.Range("H77").Value = 1120.2778
$ws.Range("I77").Value = 1004.4
$ws.Range("K77").Value = 5022
$ws.Range("M77").Value = -654

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1740.5834
$ws.Range("I20").Value = 1604.4
$ws.Range("J20").Value = 2421.5
$ws.Range("K20").Value = 1604.4
$ws.Range("L20").Value = 2421.5
$ws.Range("M20").Value = -1357.4
$ws.Range("N20").Value = -2915.5
# Row 41
$ws.Range("H41").Value = 60000
$ws.Range("J41").Value = 60000
$ws.Range("L41").Value = 60000
$ws.Range("N41").Value = -60776
# Row 47
$ws.Range("H47").Value = 60000
$ws.Range("J47").Value = 60000
$ws.Range("L47").Value = 60000
$ws.Range("N47").Value = -61040
# Row 48
$ws.Range("H48").Value = 60000
$ws.Range("J48").Value = 60000
$ws.Range("L48").Value = 60000
$ws.Range("N48").Value = -60830
# Row 86
$ws.Range("H86").Value = 2123.875
$ws.Range("I86").Value = 1998.75
$ws.Range("K86").Value = 1998.75
$ws.Range("M86").Value = -875.75
# Row 89
$ws.Range("H89").Value = 2123.875
$ws.Range("I89").Value = 1998.75
$ws.Range("K89").Value = 9993.75
$ws.Range("M89").Value = -4377.75
# Row 134
$ws.Range("H134").Value = 2302.45
$ws.Range("I134").Value = 1861.8667
$ws.Range("K134").Value = 5585.6001
$ws.Range("M134").Value = -3050.6001

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 362.22223
$ws.Range("I22").Value = 248
$ws.Range("K22").Value = 248
$ws.Range("M22").Value = 102
# Row 28
$ws.Range("H28").Value = 10274.5
$ws.Range("J28").Value = 10366
$ws.Range("L28").Value = 10366
$ws.Range("N28").Value = -10856
# Row 41
$ws.Range("H41").Value = 25086.076
$ws.Range("J41").Value = 29419.908
$ws.Range("L41").Value = 29419.908
$ws.Range("N41").Value = -30275.908
# Row 70
$ws.Range("H70").Value = 39999.332
$ws.Range("J70").Value = 39999.332
$ws.Range("L70").Value = 39999.332
$ws.Range("N70").Value = -40629.332
# Row 73
$ws.Range("H73").Value = 39999.332
$ws.Range("J73").Value = 39999.332
$ws.Range("L73").Value = 39999.332
$ws.Range("N73").Value = -42183.332
# Row 107
$ws.Range("H107").Value = 916.4
$ws.Range("I107").Value = 930.5
$ws.Range("K107").Value = 930.5
$ws.Range("M107").Value = 989.5
# Row 134
$ws.Range("H134").Value = 1685.1666
$ws.Range("I134").Value = 1496.8684
$ws.Range("K134").Value = 4490.6052
$ws.Range("M134").Value = -1955.6052

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1641.6666
$ws.Range("J34").Value = 3000
$ws.Range("L34").Value = 9000
$ws.Range("N34").Value = -9168
# Row 132
$ws.Range("H132").Value = 4107.3335
$ws.Range("I132").Value = 4951.8887
$ws.Range("J132").Value = 3262.7778
$ws.Range("K132").Value = 44566.99830000001
$ws.Range("L132").Value = 29365.0002
$ws.Range("M132").Value = -42036.99830000001
$ws.Range("N132").Value = -34425.00019999999
# Row 141
$ws.Range("H141").Value = 4805.6
$ws.Range("I141").Value = 4805.6
$ws.Range("K141").Value = 14416.8
$ws.Range("M141").Value = -9236.800000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 417
$ws.Range("I2").Value = 10.285714
$ws.Range("K2").Value = 10.285714
$ws.Range("M2").Value = 102.714286
# Row 122
$ws.Range("H122").Value = 31517.412
$ws.Range("I122").Value = 1678.0769
$ws.Range("J122").Value = 128495.25
$ws.Range("K122").Value = 5034.2307
$ws.Range("L122").Value = 385485.75
$ws.Range("M122").Value = -2584.2307
$ws.Range("N122").Value = -390385.75
# Row 134
$ws.Range("H134").Value = 110991.336
$ws.Range("J134").Value = 110991.336
$ws.Range("L134").Value = 332974.008
$ws.Range("N134").Value = -338044.008

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3704.7144
$ws.Range("I22").Value = 2840
$ws.Range("J22").Value = 5866.5
$ws.Range("K22").Value = 2840
$ws.Range("L22").Value = 5866.5
$ws.Range("M22").Value = -2545
$ws.Range("N22").Value = -6456.5
# Row 27
$ws.Range("H27").Value = 3704.7144
$ws.Range("I27").Value = 2840
$ws.Range("J27").Value = 5866.5
$ws.Range("K27").Value = 2840
$ws.Range("L27").Value = 5866.5
$ws.Range("M27").Value = -2733
$ws.Range("N27").Value = -6080.5
# Row 46
$ws.Range("I46").Value = 2916.6667
$ws.Range("J46").Value = 5297.8
$ws.Range("K46").Value = 2916.6667
$ws.Range("L46").Value = 5297.8
$ws.Range("M46").Value = -2728.6667
$ws.Range("N46").Value = -5673.8
# Row 61
$ws.Range("H61").Value = 2976.4375
$ws.Range("I61").Value = 2841.5334
$ws.Range("K61").Value = 2841.5334
$ws.Range("M61").Value = -2639.5334
# Row 113
$ws.Range("H113").Value = 2976.4375
$ws.Range("I113").Value = 2841.5334
$ws.Range("K113").Value = 2841.5334
$ws.Range("M113").Value = -671.5333999999998
# Row 132
$ws.Range("H132").Value = 5498.75
$ws.Range("J132").Value = 9995
$ws.Range("L132").Value = 29985
$ws.Range("N132").Value = -35045
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 420
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = -724
# Row 81
$ws.Range("H81").Value = 7814.5
$ws.Range("I81").Value = 7552.4287
$ws.Range("J81").Value = 8426
$ws.Range("K81").Value = 15104.8574
$ws.Range("L81").Value = 16852
$ws.Range("M81").Value = -14043.8574
$ws.Range("N81").Value = -18974
# Row 84
$ws.Range("H84").Value = 7814.5
$ws.Range("I84").Value = 7552.4287
$ws.Range("J84").Value = 8426
$ws.Range("K84").Value = 75524.28700000001
$ws.Range("L84").Value = 84260
$ws.Range("M84").Value = -70220.28700000001
$ws.Range("N84").Value = -94868
# Row 92
$ws.Range("H92").Value = 18658.5
$ws.Range("J92").Value = 18658.5
$ws.Range("L92").Value = 18658.5
$ws.Range("N92").Value = -23650.5
# Row 122
$ws.Range("H122").Value = 2318.889
$ws.Range("I122").Value = 2233.75
$ws.Range("K122").Value = 6701.25
$ws.Range("M122").Value = -4251.25

